$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the additional "2508" period, right below the existing
# data row (row 16), pushing the signature block (old rows 21-22) down to 22-23.
$ws.Rows.Item(17).Insert()

# Copy the formatting of the existing data row onto the new row so it keeps
# the same borders / fonts / fills as the rest of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new period's data (same worker, new period 2508).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73142867"
$ws.Range("D17").Value = "LUIS LEONARDO LAMBIS CAMARGO"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Update the summary figures: total "Valor Mora" now covers both periods,
# and "Cant. Periodos" goes from 1 to 2.
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
